$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 99

$ws.Cells.Item($row, 1).Value = 9
$ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($row, 3).Value = "Metropolitana"
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 4).Value = Get-Date -Year 2022 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item($row, 5).Value = 13
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101004
$ws.Cells.Item($row, 10).Value = "Frambuesa"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 380
$ws.Cells.Item($row, 14).Value = 8000
$ws.Cells.Item($row, 15).Value = 8000
$ws.Cells.Item($row, 16).Value = 8000
$ws.Cells.Item($row, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia de Linares"
$ws.Cells.Item($row, 19).Value = 4000
$ws.Cells.Item($row, 20).Value = 2
